$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1166.6666
$ws.Range("I62").Value = 1000
$ws.Range("J62").Value = 1500
$ws.Range("K62").Value = 1000
$ws.Range("L62").Value = 1500
$ws.Range("M62").Value = -376
$ws.Range("N62").Value = -2748

$ws.Range("H65").Value = 1166.6666
$ws.Range("I65").Value = 1000
$ws.Range("J65").Value = 1500
$ws.Range("K65").Value = 5000
$ws.Range("L65").Value = 7500
$ws.Range("M65").Value = -1880

$ws.Range("H135").Value = 1184.4
$ws.Range("I135").Value = 977.7143
$ws.Range("J135").Value = 1666.6666
$ws.Range("K135").Value = 8799.4287
$ws.Range("L135").Value = 14999.9994
$ws.Range("M135").Value = -6264.4287
$ws.Range("N135").Value = -20069.9994

$ws.Range("H137").Value = 776
$ws.Range("I137").Value = 776
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 2328
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 222
$ws.Range("N137").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 12431.286
$ws.Range("I38").Value = 3803.8
$ws.Range("K38").Value = 3803.8
$ws.Range("M38").Value = -3336.8

$ws.Range("H74").Value = 2148.75
$ws.Range("I74").Value = 2115.7144
$ws.Range("K74").Value = 2115.7144
$ws.Range("M74").Value = -1241.7144

$ws.Range("H77").Value = 2148.75
$ws.Range("I77").Value = 2115.7144
$ws.Range("K77").Value = 10578.572
$ws.Range("M77").Value = -6210.572

$ws.Range("H95").Value = 14208
$ws.Range("J95").Value = 14208
$ws.Range("L95").Value = 14208
$ws.Range("N95").Value = -19700

$ws.Range("H102").Value = 933.3333
$ws.Range("I102").Value = 933.3333
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 933.3333
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 688.6667
$ws.Range("N102").Value = ""

$ws.Range("H120").Value = 40000
$ws.Range("J120").Value = 40000
$ws.Range("L120").Value = 40000
$ws.Range("N120").Value = -49676

$ws.Range("H132").Value = 6881.8184
$ws.Range("I132").Value = 6881.8184
$ws.Range("K132").Value = 20645.4552
$ws.Range("M132").Value = -18115.4552

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1625.5714
$ws.Range("I86").Value = 1563.1666
$ws.Range("K86").Value = 1563.1666
$ws.Range("M86").Value = -440.1666

$ws.Range("H89").Value = 1625.5714
$ws.Range("I89").Value = 1563.1666
$ws.Range("K89").Value = 7815.833000000001
$ws.Range("M89").Value = -2199.833000000001

$ws.Range("H94").Value = 1724.4642
$ws.Range("I94").Value = 1422.1904
$ws.Range("J94").Value = 2631.2856
$ws.Range("K94").Value = 1422.1904
$ws.Range("L94").Value = 2631.2856
$ws.Range("M94").Value = -971.1904
$ws.Range("N94").Value = -3533.2856

$ws.Range("H105").Value = 2483.3333
$ws.Range("J105").Value = 2000
$ws.Range("L105").Value = 2000
$ws.Range("N105").Value = -5494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1265.6666
$ws.Range("I31").Value = 1265.6666
$ws.Range("K31").Value = 1265.6666
$ws.Range("M31").Value = -970.6666

$ws.Range("H34").Value = 1265.6666
$ws.Range("I34").Value = 1265.6666
$ws.Range("K34").Value = 1265.6666
$ws.Range("M34").Value = -1063.6666

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = ""

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = ""

$ws.Range("H99").Value = 836082.7
$ws.Range("I99").Value = 2874
$ws.Range("K99").Value = 2874
$ws.Range("M99").Value = -1376

$ws.Range("H122").Value = 6114.8335
$ws.Range("I122").Value = 2247.5
$ws.Range("J122").Value = 8048.5
$ws.Range("K122").Value = 6742.5
$ws.Range("L122").Value = 24145.5
$ws.Range("M122").Value = -4292.5
$ws.Range("N122").Value = -29045.5

$ws.Range("H126").Value = 836082.7
$ws.Range("I126").Value = 2874
$ws.Range("K126").Value = 8622
$ws.Range("M126").Value = -6152

$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").Value = ""

$ws.Range("H134").Value = 2529.3333
$ws.Range("I134").Value = 2529.3333
$ws.Range("K134").Value = 7587.999899999999
$ws.Range("M134").Value = -5052.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 246.88889
$ws.Range("I7").Value = 215
$ws.Range("J7").Value = 502
$ws.Range("K7").Value = 645
$ws.Range("L7").Value = 1506
$ws.Range("M7").Value = -533
$ws.Range("N7").Value = -1730

$ws.Range("H35").Value = 16829.166
$ws.Range("I35").Value = 487.5
$ws.Range("K35").Value = 1462.5
$ws.Range("M35").Value = -1174.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2841.7
$ws.Range("I122").Value = 2070
$ws.Range("J122").Value = 3999.25
$ws.Range("K122").Value = 6210
$ws.Range("L122").Value = 11997.75
$ws.Range("M122").Value = -3760
$ws.Range("N122").Value = -16897.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 111
$ws.Range("N2").Value = ""

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = ""

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = ""

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = ""

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = ""

$ws.Range("H133").Value = 120000
$ws.Range("J133").Value = 120000
$ws.Range("L133").Value = 120000
$ws.Range("N133").Value = -125060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 4016357.2
$ws.Range("I5").Value = 5550000.5
$ws.Range("J5").Value = 3402900
$ws.Range("K5").Value = 5550000.5
$ws.Range("L5").Value = 3402900
$ws.Range("M5").Value = -5549888.5
$ws.Range("N5").Value = -3403124

$ws.Range("H19").Value = 16006
$ws.Range("J19").Value = 16006
$ws.Range("L19").Value = 16006
$ws.Range("N19").Value = -16354

$ws.Range("H29").Value = 9999.5
$ws.Range("I29").Value = 9999.5
$ws.Range("K29").Value = 9999.5
$ws.Range("M29").Value = -9709.5

$ws.Range("H122").Value = 1805.3334
$ws.Range("I122").Value = 1710.3334
$ws.Range("J122").Value = 1995.3334
$ws.Range("K122").Value = 5131.0002
$ws.Range("L122").Value = 5986.0002
$ws.Range("M122").Value = -2681.0002
$ws.Range("N122").Value = -10886.0002

$ws.Range("H132").Value = 1774
$ws.Range("I132").Value = 1956
$ws.Range("K132").Value = 5868
$ws.Range("M132").Value = -3338
